{"js": "// 1. Merge the title runs \"Coaster Capacity Calculator\" + \" (\" + \"Final Assignment\" + \")\"\n//    into a single run reading \"Coaster Capacity Calculator (Final Assignment)\".\nconst body = context.document.body;\nconst titleHits = body.search(\"Coaster Capacity Calculator (Final Assignment)\", { matchCase: true, matchWildcards: false });\ntitleHits.load(\"items\");\nawait context.sync();\n\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText(\"Coaster Capacity Calculator (Final Assignment)\", \"Replace\");\n  await context.sync();\n}\n\n// 2. Add a new table row for \"radSpecRestraint\" right after the \"radOTS\" row.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  row.load(\"values\");\n}\nawait context.sync();\n\nlet radOTSRow = null;\nfor (const row of table.rows.items) {\n  const firstCellText = row.values[0][0];\n  if (firstCellText === \"radOTS\") {\n    radOTSRow = row;\n    break;\n  }\n}\n\nif (radOTSRow) {\n  radOTSRow.insertRows(\"After\", 1, [[\"radSpecRestraint\", \"none\", \"operated externally\"]]);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Merge the title runs: \"Coaster Capacity Calculator\" + \" (\" + \"Final Assignment\" + \")\"\n#        into a single run reading \"Coaster Capacity Calculator (Final Assignment)\".\n$find = $d.Range().Find\n$find.Text = \"Coaster Capacity Calculator (Final Assignment)\"\n$find.MatchWildcards = $false\n$find.Replacement.Text = \"Coaster Capacity Calculator (Final Assignment)\"\n# wdReplaceAll = 2\n$find.Execute($null,$false,$false,$false,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n# --- 2. Add a new table row for \"radSpecRestraint\" right after the \"radOTS\" row.\n$t = $d.Tables.Item(1)\n$radOTSRow = $null\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    if ($t.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13) -eq \"radOTS\") {\n        $radOTSRow = $i\n        break\n    }\n}\n\n$newRow = $t.Rows.Add($t.Rows.Item($radOTSRow + 1))\n$newRow.Cells.Item(1).Range.Text = \"radSpecRestraint\"\n$newRow.Cells.Item(2).Range.Text = \"none\"\n$newRow.Cells.Item(3).Range.Text = \"operated externally\"\n"}
